$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete columns G then F (rightmost two columns) to drop the
# "chan-xiao-lv" (F) and "xiao-shou-liang" (G) columns entirely, shrinking
# the used range from A1:G81 to A1:E81.
$ws.Columns.Item(7).Delete()
$ws.Columns.Item(6).Delete()

# Step 2: within every 4-row year block (A/B/C/D sub-periods), the "B"
# sub-period row and the "C" sub-period row have their columns A:E content
# swapped (row numbers stay put; only the data moves).
# Year block starting row 2: swap row 3 and row 4
$ws.Cells.Item(3,1).Value = "2000年C"
$ws.Cells.Item(4,1).Value = "2000年B"
$ws.Cells.Item(3,2).Value = 99.09999999999999
$ws.Cells.Item(4,2).Value = 98.5
$ws.Cells.Item(3,4).Value = 3.2
$ws.Cells.Item(4,4).Value = 3
$ws.Cells.Item(3,5).Value = 309.8
$ws.Cells.Item(4,5).Value = 216.7

# Year block starting row 6: swap row 7 and row 8
$ws.Cells.Item(7,1).Value = "2001年C"
$ws.Cells.Item(8,1).Value = "2001年B"
$ws.Cells.Item(7,2).Value = 101.5
$ws.Cells.Item(8,2).Value = 100.5
$ws.Cells.Item(7,3).Value = 2.4
$ws.Cells.Item(8,3).Value = ""
$ws.Cells.Item(7,4).Value = -3.4
$ws.Cells.Item(8,4).Value = -1.1
$ws.Cells.Item(7,5).Value = 295.3
$ws.Cells.Item(8,5).Value = 212.9

# Year block starting row 10: swap row 11 and row 12
$ws.Cells.Item(11,1).Value = "2002年C"
$ws.Cells.Item(12,1).Value = "2002年B"
$ws.Cells.Item(11,2).Value = 101.8
$ws.Cells.Item(12,2).Value = 100
$ws.Cells.Item(11,4).Value = -5.7
$ws.Cells.Item(12,4).Value = 0.1
$ws.Cells.Item(11,5).Value = 259.6
$ws.Cells.Item(12,5).Value = 179.7

# Year block starting row 14: swap row 15 and row 16
$ws.Cells.Item(15,1).Value = "2003年C"
$ws.Cells.Item(16,1).Value = "2003年B"
$ws.Cells.Item(15,2).Value = 101.7
$ws.Cells.Item(16,2).Value = 101.4
$ws.Cells.Item(15,4).Value = -5.4
$ws.Cells.Item(16,4).Value = -3.8
$ws.Cells.Item(15,5).Value = 226.7
$ws.Cells.Item(16,5).Value = 153.4

# Year block starting row 18: swap row 19 and row 20
$ws.Cells.Item(19,1).Value = "2004年C"
$ws.Cells.Item(20,1).Value = "2004年B"
$ws.Cells.Item(19,2).Value = 101.9
$ws.Cells.Item(20,2).Value = 102
$ws.Cells.Item(19,4).Value = -2.9
$ws.Cells.Item(20,4).Value = -4.9
$ws.Cells.Item(19,5).Value = 215.2
$ws.Cells.Item(20,5).Value = 138.9

# Year block starting row 22: swap row 23 and row 24
$ws.Cells.Item(23,1).Value = "2005年C"
$ws.Cells.Item(24,1).Value = "2005年B"
$ws.Cells.Item(23,2).Value = 98.40000000000001
$ws.Cells.Item(24,2).Value = 100.8
$ws.Cells.Item(23,3).Value = -3.5
$ws.Cells.Item(24,3).Value = -1.2
$ws.Cells.Item(23,4).Value = -0.9
$ws.Cells.Item(24,4).Value = -2.6
$ws.Cells.Item(23,5).Value = 230.7
$ws.Cells.Item(24,5).Value = 154.7

# Year block starting row 26: swap row 27 and row 28
$ws.Cells.Item(27,1).Value = "2006年C"
$ws.Cells.Item(28,1).Value = "2006年B"
$ws.Cells.Item(27,2).Value = 97.3
$ws.Cells.Item(28,2).Value = 97.90000000000001
$ws.Cells.Item(27,3).Value = -1.1
$ws.Cells.Item(28,3).Value = -2.9
$ws.Cells.Item(27,4).Value = 8
$ws.Cells.Item(28,4).Value = 5.7
$ws.Cells.Item(27,5).Value = 264.9
$ws.Cells.Item(28,5).Value = 178.3

# Year block starting row 30: swap row 31 and row 32
$ws.Cells.Item(31,1).Value = "2007年C"
$ws.Cells.Item(32,1).Value = "2007年B"
$ws.Cells.Item(31,2).Value = 98.59999999999999
$ws.Cells.Item(32,2).Value = 98.8
$ws.Cells.Item(31,3).Value = 1.3
$ws.Cells.Item(32,3).Value = 0.9
$ws.Cells.Item(31,4).Value = 3.6
$ws.Cells.Item(32,4).Value = 1.6
$ws.Cells.Item(31,5).Value = 329
$ws.Cells.Item(32,5).Value = 211.8

# Year block starting row 34: swap row 35 and row 36
$ws.Cells.Item(35,1).Value = "2008年C"
$ws.Cells.Item(36,1).Value = "2008年B"
$ws.Cells.Item(35,2).Value = 97.59999999999999
$ws.Cells.Item(36,2).Value = 97.5
$ws.Cells.Item(35,3).Value = -0.7
$ws.Cells.Item(36,3).Value = -1.1
$ws.Cells.Item(35,4).Value = 12.1
$ws.Cells.Item(36,4).Value = 5.8
$ws.Cells.Item(35,5).Value = 378
$ws.Cells.Item(36,5).Value = 252.7

# Year block starting row 38: swap row 39 and row 40
$ws.Cells.Item(39,1).Value = "2009年C"
$ws.Cells.Item(40,1).Value = "2009年B"
$ws.Cells.Item(39,2).Value = 99.2
$ws.Cells.Item(40,2).Value = 99.90000000000001
$ws.Cells.Item(39,3).Value = -2.6
$ws.Cells.Item(40,3).Value = 2.7
$ws.Cells.Item(39,4).Value = 9.1
$ws.Cells.Item(40,4).Value = 1.8
$ws.Cells.Item(39,5).Value = 434.7
$ws.Cells.Item(40,5).Value = 309.3

# Year block starting row 42: swap row 43 and row 44
$ws.Cells.Item(43,1).Value = "2010年C"
$ws.Cells.Item(44,1).Value = "2010年B"
$ws.Cells.Item(43,2).Value = 98.2
$ws.Cells.Item(44,2).Value = 99.5
$ws.Cells.Item(43,3).Value = -2.6
$ws.Cells.Item(44,3).Value = -0.5
$ws.Cells.Item(43,4).Value = 8.300000000000001
$ws.Cells.Item(44,4).Value = 0.2
$ws.Cells.Item(43,5).Value = 597.9
$ws.Cells.Item(44,5).Value = 408.4

# Year block starting row 46: swap row 47 and row 48
$ws.Cells.Item(47,1).Value = "2011年C"
$ws.Cells.Item(48,1).Value = "2011年B"
$ws.Cells.Item(47,2).Value = 97
$ws.Cells.Item(48,2).Value = 98.90000000000001
$ws.Cells.Item(47,3).Value = -0.4
$ws.Cells.Item(48,3).Value = 0.2
$ws.Cells.Item(47,4).Value = 12.4
$ws.Cells.Item(48,4).Value = 6.5
$ws.Cells.Item(47,5).Value = 689
$ws.Cells.Item(48,5).Value = 471.9

# Year block starting row 50: swap row 51 and row 52
$ws.Cells.Item(51,1).Value = "2012年C"
$ws.Cells.Item(52,1).Value = "2012年B"
$ws.Cells.Item(51,2).Value = 99.90000000000001
$ws.Cells.Item(52,2).Value = 99
$ws.Cells.Item(51,3).Value = 3.3
$ws.Cells.Item(52,3).Value = 1.7
$ws.Cells.Item(51,4).Value = 2.9
$ws.Cells.Item(52,4).Value = 8.300000000000001
$ws.Cells.Item(51,5).Value = 794.5
$ws.Cells.Item(52,5).Value = 528.1

# Year block starting row 54: swap row 55 and row 56
$ws.Cells.Item(55,1).Value = "2013年C"
$ws.Cells.Item(56,1).Value = "2013年B"
$ws.Cells.Item(55,2).Value = 96.5
$ws.Cells.Item(56,2).Value = 97.09999999999999
$ws.Cells.Item(55,3).Value = -2.5
$ws.Cells.Item(56,3).Value = -0.9
$ws.Cells.Item(55,4).Value = 24
$ws.Cells.Item(56,4).Value = 13.6
$ws.Cells.Item(55,5).Value = 833.6
$ws.Cells.Item(56,5).Value = 569.2

# Year block starting row 58: swap row 59 and row 60
$ws.Cells.Item(59,1).Value = "2014年C"
$ws.Cells.Item(60,1).Value = "2014年B"
$ws.Cells.Item(59,2).Value = 96.09999999999999
$ws.Cells.Item(60,2).Value = 96.8
$ws.Cells.Item(59,3).Value = -0.2
$ws.Cells.Item(60,3).Value = 0.2
$ws.Cells.Item(59,4).Value = 20.1
$ws.Cells.Item(60,4).Value = 10.6
$ws.Cells.Item(59,5).Value = 855.5
$ws.Cells.Item(60,5).Value = 591

# Year block starting row 62: swap row 63 and row 64
$ws.Cells.Item(63,1).Value = "2015年C"
$ws.Cells.Item(64,1).Value = "2015年B"
$ws.Cells.Item(63,2).Value = 97.90000000000001
$ws.Cells.Item(64,2).Value = 98.7
$ws.Cells.Item(63,3).Value = 2.5
$ws.Cells.Item(64,3).Value = 2.3
$ws.Cells.Item(63,4).Value = 4.4
$ws.Cells.Item(64,4).Value = 2.6
$ws.Cells.Item(63,5).Value = 914.4
$ws.Cells.Item(64,5).Value = 611.5

# Year block starting row 66: swap row 67 and row 68
$ws.Cells.Item(67,1).Value = "2016年C"
$ws.Cells.Item(68,1).Value = "2016年B"
$ws.Cells.Item(67,2).Value = 97.09999999999999
$ws.Cells.Item(68,2).Value = 98.90000000000001
$ws.Cells.Item(67,3).Value = -0.6
$ws.Cells.Item(68,3).Value = -0.3
$ws.Cells.Item(67,4).Value = 9.9
$ws.Cells.Item(68,4).Value = 3.3
$ws.Cells.Item(67,5).Value = 939.01257
$ws.Cells.Item(68,5).Value = 643.57245

# Year block starting row 70: swap row 71 and row 72
$ws.Cells.Item(71,1).Value = "2017年C"
$ws.Cells.Item(72,1).Value = "2017年B"
$ws.Cells.Item(71,2).Value = 99.09999999999999
$ws.Cells.Item(72,2).Value = 101.1
$ws.Cells.Item(71,3).Value = 1.7
$ws.Cells.Item(72,3).Value = 2.6
$ws.Cells.Item(71,4).Value = 1.6
$ws.Cells.Item(72,4).Value = -3.2
$ws.Cells.Item(71,5).Value = 979.0948
$ws.Cells.Item(72,5).Value = 685.30074

# Year block starting row 74: swap row 75 and row 76
$ws.Cells.Item(75,1).Value = "2018年C"
$ws.Cells.Item(76,1).Value = "2018年B"
$ws.Cells.Item(75,2).Value = 103.8
$ws.Cells.Item(76,2).Value = 103
$ws.Cells.Item(75,3).Value = 4.1
$ws.Cells.Item(76,3).Value = -1.3
$ws.Cells.Item(75,4).Value = -2.3
$ws.Cells.Item(76,4).Value = -6.6
$ws.Cells.Item(75,5).Value = 700.82303
$ws.Cells.Item(76,5).Value = 501.85575

# Year block starting row 78: swap row 79 and row 80
$ws.Cells.Item(79,1).Value = "2019年C"
$ws.Cells.Item(80,1).Value = "2019年B"
$ws.Cells.Item(79,2).Value = 98.59999999999999
$ws.Cells.Item(80,2).Value = 101.1
$ws.Cells.Item(79,3).Value = -0.1
$ws.Cells.Item(80,3).Value = 2.3
$ws.Cells.Item(79,4).Value = 3.7
$ws.Cells.Item(80,4).Value = -2.8
$ws.Cells.Item(79,5).Value = 561.17033
$ws.Cells.Item(80,5).Value = 397.96866

"done"